# Author: "Have learned more topics"
# Adds a new "Data Structures" section header above the existing topic
# list, and a brand-new "Algorithms" section (with its own header +
# topics) below it, following the same visual formatting already used
# in the sheet (title/header/data-row styles, banding borders, etc).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Insert a new row above the "Topic / Status" header row (old row 2)
#    to host the "Data Structures" sub-title, pushing everything below
#    down by one row.
# ---------------------------------------------------------------------
$ws.Rows.Item(2).Insert()

# Fill in the new row 2 content + merge, matching the look of the big
# "DSA" title in row 1 but as a black-on-gold section banner.
$ws.Range("A2:B2").Merge()
$ws.Range("A2").Value2 = "Data Structures"

$secRange = $ws.Range("A2:B2")
$secRange.Font.Bold = $true
$secRange.Font.ThemeColor = 1
$secRange.Interior.ThemeColor = 8
$secRange.HorizontalAlignment = -4108

# ---------------------------------------------------------------------
# 2) Row 1 ("DSA" banner) loses its bottom border (it now sits directly
#    above the new "Data Structures" banner), row 3 ("Topic /
#    Status" header) loses its top border (it now sits directly below
#    the borderless "Data Structures" banner).
# ---------------------------------------------------------------------
foreach ($addr in @("A1", "B1")) {
    $c = $ws.Range($addr)
    $c.Borders.Item(9).LineStyle = -4142
}
foreach ($addr in @("A3", "B3")) {
    $c = $ws.Range($addr)
    $c.Borders.Item(8).LineStyle = -4142
}

# ---------------------------------------------------------------------
# 3) Append two blank data-styled rows (16 & 17) after "Tree", matching
#    the plain data-row look already used for rows 4-15.
# ---------------------------------------------------------------------
$ws.Range("A15:B15").Copy()
$ws.Range("A16:B17").PasteSpecial(-4122)
$ws.Range("A16").ClearContents()
$ws.Range("A17").ClearContents()

# ---------------------------------------------------------------------
# 4) Row 18: a blank "cap" row that only has a top border (left/right/
#    top thin, no bottom) - same data-row borderless background but the
#    top of a new bordered box starting the "Algorithms" section.
# ---------------------------------------------------------------------
$ws.Range("A16:B16").Copy()
$ws.Range("A18:B18").PasteSpecial(-4122)
foreach ($addr in @("A18", "B18")) {
    $c = $ws.Range($addr)
    $c.Borders.Item(7).LineStyle = 1
    $c.Borders.Item(7).Weight = 2
    $c.Borders.Item(10).LineStyle = 1
    $c.Borders.Item(10).Weight = 2
    $c.Borders.Item(8).LineStyle = 1
    $c.Borders.Item(8).Weight = 2
    $c.Borders.Item(9).LineStyle = -4142
}

# ---------------------------------------------------------------------
# 5) Row 19: "Algorithms" section banner - same dark teal "title" look
#    as row 1 / the old title row, merged + centered, no border.
# ---------------------------------------------------------------------
$ws.Range("A1:B1").Copy()
$ws.Range("A19:B19").PasteSpecial(-4122)
$ws.Range("A19:B19").Merge()
$ws.Range("A19").Value2 = "Algorithms"
foreach ($addr in @("A19", "B19")) {
    $c = $ws.Range($addr)
    $c.Borders.Item(7).LineStyle = -4142
    $c.Borders.Item(10).LineStyle = -4142
    $c.Borders.Item(8).LineStyle = -4142
    $c.Borders.Item(9).LineStyle = -4142
}

# ---------------------------------------------------------------------
# 6) Row 20: first Algorithms topic ("Linear Search") - plain data row
#    but capped with only a bottom border (no top), closing the little
#    box opened at row 18.
# ---------------------------------------------------------------------
$ws.Range("A4:B4").Copy()
$ws.Range("A20:B20").PasteSpecial(-4122)
$ws.Range("A20").Value2 = "Linear Search"
$ws.Range("B20").Value2 = "Done"
foreach ($addr in @("A20", "B20")) {
    $c = $ws.Range($addr)
    $c.Borders.Item(8).LineStyle = -4142
}

# ---------------------------------------------------------------------
# 7) Rows 21-24: remaining Algorithms topics, same plain data-row style
#    used throughout the rest of the sheet.
# ---------------------------------------------------------------------
$ws.Range("A4:B4").Copy()
$ws.Range("A21:B24").PasteSpecial(-4122)

$ws.Range("A21").Value2 = "Binary Search"
$ws.Range("B21").Value2 = "Done"
$ws.Range("A22").Value2 = "Bubble Sort"
$ws.Range("B22").Value2 = "Done"
$ws.Range("A23").Value2 = "Selection Sort"
$ws.Range("B23").Value2 = "Done"
$ws.Range("A24").Value2 = "Merge Sort"

# ---------------------------------------------------------------------
# 8) Row 25: trailing blank data row, same as rows 16/17 above.
# ---------------------------------------------------------------------
$ws.Range("A4:B4").Copy()
$ws.Range("A25:B25").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 9) Restore the view / selection to match the authored state.
# ---------------------------------------------------------------------
$ws.Range("C25").Select()
$ws.Application.ActiveWindow.ScrollRow = 23
